# Update the Fitness values (column C) in the active worksheet.
# Column A = Run, Column B = Generation, Column C = Fitness
# The new fitness values follow a step pattern across the generation rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ranges = @(
    @{ Start = 2;   End = 3;   Value = 11078 },
    @{ Start = 4;   End = 5;   Value = 10381 },
    @{ Start = 6;   End = 11;  Value = 9852 },
    @{ Start = 12;  End = 21;  Value = 8996 },
    @{ Start = 22;  End = 26;  Value = 8938 },
    @{ Start = 27;  End = 40;  Value = 8446 },
    @{ Start = 41;  End = 42;  Value = 7884 },
    @{ Start = 43;  End = 46;  Value = 7882 },
    @{ Start = 47;  End = 252; Value = 7310 }
)

foreach ($r in $ranges) {
    $ws.Range("C$($r.Start):C$($r.End)").Value = $r.Value
}
